# CMSIS-Driver USART retarget implementation: add a "standard stream" entry
# to the "Custom" retarget options list on slide 3 of the block diagram,
# renaming the previous "Custom" entry to "UART" and growing the containing
# rectangles to make room for the new entry.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$EMU_PER_POINT = 12700.0

# 1) Grow the vertical "Retarget" label rectangle on the left diagram so it
#    still spans the full height of the (now taller) stack of boxes below it.
$retargetLabel = $s.Shapes.Item("Rectangle 1")   # id=2
$retargetLabel.Height = 4096687 / $EMU_PER_POINT

# 2) Grow the dark-blue container rectangle that holds the list of retarget
#    backends (Breakpoint / Event Recorder / ITM / Custom) so it can fit the
#    additional "Custom" entry being appended below.
$container = $s.Shapes.Item("Rectangle 6")       # id=7
$container.Height = 2699028 / $EMU_PER_POINT

# 3) Rename the existing "Custom" backend entry to "UART".
$uartBox = $s.Shapes.Item("Rectangle 11")        # id=12, text "Custom"
$uartBox.TextFrame.TextRange.Text = "UART"

# 4) Add a new "Custom" backend entry below the others, matching the style
#    of the existing boxes (duplicate the equivalent box from the right-hand
#    diagram, then reposition/rename it).
$template = $s.Shapes.Item("Rectangle 30")       # id=31, text "Custom"
$newRange = $template.Duplicate()
$newBox = $newRange.Item(1)
$newBox.Name = "Rectangle 5"
$newBox.Left = 1194700 / $EMU_PER_POINT
$newBox.Top = 3645899 / $EMU_PER_POINT
# Width/Height are already correct (1615439 x 364464 EMU) via Duplicate(),
# so leave them untouched to avoid point-conversion rounding drift.
$newBox.TextFrame.TextRange.Text = "Custom"
